$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The codeforiati:group-code / codeforiati:group-name columns (D and E)
# had their contents swapped: column D now holds what used to be in
# column E, and column E now holds what used to be in column D - for the
# header row and every data row in the sheet's used range.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $eVal
    $ws.Cells.Item($r, 5).Value2 = $dVal
}
